$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.617.53'
$ws.Range("E2").Value = '  -5.02%  '
$ws.Range("D3").Value = '3.013.07'
$ws.Range("E3").Value = '  -6.30%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.66'
$ws.Range("E5").Value = '  -2.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.50'
$ws.Range("E6").Value = '  -7.22%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").Value = '3.003.78'
$ws.Range("E8").Value = '  -6.58%  '
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.133'
$ws.Range("E10").Value = '  -7.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.16'
$ws.Range("E11").Value = '  -2.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  -3.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000224'
$ws.Range("E13").Value = '  -7.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.76'
$ws.Range("E14").Value = '  -6.74%  '
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").Value = '3.509.63'
$ws.Range("E16").Value = '  -6.27%  '
$ws.Range("D17").Value = '3.012.78'
$ws.Range("E17").Value = '  -6.25%  '
$ws.Range("D18").Value = '60.538.26'
$ws.Range("E18").Value = '  -5.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.43'
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '433.62'
$ws.Range("E20").Value = '  -7.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.16'
$ws.Range("E21").Value = '  -6.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.666'
$ws.Range("E22").Value = '  -5.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.04'
$ws.Range("E23").Value = '  -8.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.88'
$ws.Range("E24").Value = '  -5.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.57'
$ws.Range("E25").Value = '  -4.61%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.56'
$ws.Range("E28").Value = '  -5.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("E29").Value = '  -6.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.92'
$ws.Range("E30").Value = '  -7.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.20'
$ws.Range("E31").Value = '  -10.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.40'
$ws.Range("E32").Value = '  -7.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0941'
$ws.Range("E33").Value = '  -9.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.17'
$ws.Range("E34").Value = '  -10.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.962'
$ws.Range("E35").Value = '  -7.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.62'
$ws.Range("E36").Value = '  -4.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '50.20'
$ws.Range("E37").Value = '  -2.93%  '
$ws.Range("D38").Value = '0.0₃0671'
$ws.Range("E38").Value = '  -8.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.43'
$ws.Range("E39").Value = '  +3.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0362'
$ws.Range("E40").Value = '  -8.23%  '
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '387.41'
$ws.Range("E42").Value = '  -5.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.53'
$ws.Range("E43").Value = '  -8.86%  '
$ws.Range("D44").Value = '2.663.96'
$ws.Range("E44").Value = '  -5.83%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.237'
$ws.Range("E46").Value = '  -7.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.03'
$ws.Range("E47").Value = '  -6.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.17'
$ws.Range("E48").Value = '  -7.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.107'
$ws.Range("E49").Value = '  -4.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.92'
$ws.Range("E50").Value = '  -7.60%  '
$ws.Range("E51").Value = '  +3.75%  '
